$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 98 (shifts rows 98:109 down to 99:110)
$ws.Rows.Item(98).Insert()

# Fill the new row 98 with data matching the surrounding rows' pattern
$ws.Cells.Item(98, 1).Value = 8
$ws.Cells.Item(98, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(98, 3).Value = "Coquimbo"
$ws.Cells.Item(98, 4).Value = 44449
$ws.Cells.Item(98, 5).Value = 4
$ws.Cells.Item(98, 6).Value = 100112037
$ws.Cells.Item(98, 7).Value = "Cebollín"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 3080
$ws.Cells.Item(98, 11).Value = 900
$ws.Cells.Item(98, 12).Value = 1000
$ws.Cells.Item(98, 13).Value = 950
$ws.Cells.Item(98, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(98, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(98, 16).Value = 158
$ws.Cells.Item(98, 17).Value = 6
$ws.Cells.Item(98, 18).Value = "Hortaliza"

# Match the date number format/style used by the other rows in column D
$ws.Cells.Item(98, 4).NumberFormat = $ws.Cells.Item(99, 4).NumberFormat
